# "complete Branch in struction.xlsx"
#
# Fills in the "Branch" column (K) values for the remaining instruction
# rows (25-30) on Sheet1, and updates the active selection / scroll
# position left behind by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the Branch (column K) values for rows 25-30 ---------------
# R-type / jump-ish rows (no branch) keep "00" (already present, untouched).
# These six rows get their Branch control-signal value completed:
$ws.Range("K25").Value = "01"
$ws.Range("K26").Value = "01"
$ws.Range("K27").Value = "01"
$ws.Range("K28").Value = "01"
$ws.Range("K29").Value = "10"
$ws.Range("K30").Value = "10"

# --- Update view state left after editing -------------------------------
# Move the frozen-pane scroll position down toward the bottom of the table
# and leave the final selection on the last edited area of the sheet.
$win = $wb.Windows.Item(1)
$win.Left = 5820
$win.Top = 105
$win.ScrollRow = 8
$win.ScrollColumn = 7

$ws.Range("I31").Select() | Out-Null
